# Auto-generated edit script applying the Garuda_Profits.xlsx leve-profit refresh
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(137, 8).Value = 1292.7646   # H137: 1204.0526 -> 1292.7646
$ws.Cells.Item(137, 9).Value = 1198.9333   # I137: 1155.25 -> 1198.9333
$ws.Cells.Item(137, 10).Value = 1996.5   # J137: 1464.3334 -> 1996.5
$ws.Cells.Item(137, 11).Value = 3596.7999   # K137: 3465.75 -> 3596.7999
$ws.Cells.Item(137, 12).Value = 5989.5   # L137: 4393.0002 -> 5989.5
$ws.Cells.Item(137, 13).Value = -1046.7999   # M137: -915.75 -> -1046.7999
$ws.Cells.Item(137, 14).Value = -11089.5   # N137: -9493.0002 -> -11089.5
$ws.Cells.Item(138, 8).Value = 1372.683   # H138: 1409.2273 -> 1372.683
$ws.Cells.Item(138, 9).Value = 725.38464   # I138: 723 -> 725.38464
$ws.Cells.Item(138, 10).Value = 1673.2142   # J138: 2095.4546 -> 1673.2142
$ws.Cells.Item(138, 11).Value = 2176.15392   # K138: 2169 -> 2176.15392
$ws.Cells.Item(138, 12).Value = 5019.642599999999   # L138: 6286.3638 -> 5019.642599999999
$ws.Cells.Item(138, 13).Value = 2963.84608   # M138: 2971 -> 2963.84608
$ws.Cells.Item(138, 14).Value = -15299.6426   # N138: -16566.3638 -> -15299.6426
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(63, 8).Value = 1430541.2   # H63: 2001873.8 -> 1430541.2
$ws.Cells.Item(63, 9).Value = 1430541.2   # I63: 2001873.8 -> 1430541.2
$ws.Cells.Item(63, 11).Value = 1430541.2   # K63: 2001873.8 -> 1430541.2
$ws.Cells.Item(63, 13).Value = -1429855.2   # M63: -2001187.8 -> -1429855.2
$ws.Cells.Item(66, 8).Value = 1430541.2   # H66: 2001873.8 -> 1430541.2
$ws.Cells.Item(66, 9).Value = 1430541.2   # I66: 2001873.8 -> 1430541.2
$ws.Cells.Item(66, 11).Value = 7152706   # K66: 10009369 -> 7152706
$ws.Cells.Item(66, 13).Value = -7149274   # M66: -10005937 -> -7149274
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 3633.3333   # H86: 2800.5 -> 3633.3333
$ws.Cells.Item(86, 9).Value = 2950   # I86: 2800.5 -> 2950
$ws.Cells.Item(86, 10).Value = 5000   # J86: 0 -> 5000
$ws.Cells.Item(86, 11).Value = 2950   # K86: 2800.5 -> 2950
$ws.Cells.Item(86, 12).Value = 5000   # L86: 0 -> 5000
$ws.Cells.Item(86, 13).Value = -1827   # M86: -1677.5 -> -1827
$ws.Cells.Item(86, 14).Value = -7246   # N86: None -> -7246
$ws.Cells.Item(89, 8).Value = 3633.3333   # H89: 2800.5 -> 3633.3333
$ws.Cells.Item(89, 9).Value = 2950   # I89: 2800.5 -> 2950
$ws.Cells.Item(89, 10).Value = 5000   # J89: 0 -> 5000
$ws.Cells.Item(89, 11).Value = 14750   # K89: 14002.5 -> 14750
$ws.Cells.Item(89, 12).Value = 25000   # L89: 0 -> 25000
$ws.Cells.Item(89, 13).Value = -9134   # M89: -8386.5 -> -9134
$ws.Cells.Item(89, 14).Value = -36232   # N89: None -> -36232
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 1552.7084   # H58: 1960 -> 1552.7084
$ws.Cells.Item(58, 9).Value = 1148.0555   # I58: 1725.5555 -> 1148.0555
$ws.Cells.Item(58, 10).Value = 2766.6667   # J58: 2194.4443 -> 2766.6667
$ws.Cells.Item(58, 11).Value = 1148.0555   # K58: 1725.5555 -> 1148.0555
$ws.Cells.Item(58, 12).Value = 2766.6667   # L58: 2194.4443 -> 2766.6667
$ws.Cells.Item(58, 13).Value = -945.0554999999999   # M58: -1522.5555 -> -945.0554999999999
$ws.Cells.Item(58, 14).Value = -3172.6667   # N58: -2600.4443 -> -3172.6667
$ws.Cells.Item(132, 8).Value = 2283.3103   # H132: 1979.697 -> 2283.3103
$ws.Cells.Item(132, 9).Value = 1246.3182   # I132: 1127.1482 -> 1246.3182
$ws.Cells.Item(132, 10).Value = 5542.4287   # J132: 5816.1665 -> 5542.4287
$ws.Cells.Item(132, 11).Value = 3738.9546   # K132: 3381.4446 -> 3738.9546
$ws.Cells.Item(132, 12).Value = 16627.2861   # L132: 17448.4995 -> 16627.2861
$ws.Cells.Item(132, 13).Value = -1208.9546   # M132: -851.4446000000003 -> -1208.9546
$ws.Cells.Item(132, 14).Value = -21687.2861   # N132: -22508.4995 -> -21687.2861
$ws.Cells.Item(134, 8).Value = 1254.15   # H134: 1174.2916 -> 1254.15
$ws.Cells.Item(134, 9).Value = 1226.8334   # I134: 1164.6666 -> 1226.8334
$ws.Cells.Item(134, 10).Value = 1500   # J134: 1190.3334 -> 1500
$ws.Cells.Item(134, 11).Value = 3680.5002   # K134: 3493.9998 -> 3680.5002
$ws.Cells.Item(134, 12).Value = 4500   # L134: 3571.0002 -> 4500
$ws.Cells.Item(134, 13).Value = -1145.5002   # M134: -958.9998000000001 -> -1145.5002
$ws.Cells.Item(134, 14).Value = -9570   # N134: -8641.0002 -> -9570
$ws.Cells.Item(136, 8).Value = 1552.7084   # H136: 1960 -> 1552.7084
$ws.Cells.Item(136, 9).Value = 1148.0555   # I136: 1725.5555 -> 1148.0555
$ws.Cells.Item(136, 10).Value = 2766.6667   # J136: 2194.4443 -> 2766.6667
$ws.Cells.Item(136, 11).Value = 3444.1665   # K136: 5176.666499999999 -> 3444.1665
$ws.Cells.Item(136, 12).Value = 8300.000100000001   # L136: 6583.3329 -> 8300.000100000001
$ws.Cells.Item(136, 13).Value = -894.1664999999998   # M136: -2626.666499999999 -> -894.1664999999998
$ws.Cells.Item(136, 14).Value = -13400.0001   # N136: -11683.3329 -> -13400.0001
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(31, 8).Value = 600   # H31: 0 -> 600
$ws.Cells.Item(31, 9).Value = 500   # I31: 0 -> 500
$ws.Cells.Item(31, 10).Value = 700   # J31: 0 -> 700
$ws.Cells.Item(31, 11).Value = 1500   # K31: 0 -> 1500
$ws.Cells.Item(31, 12).Value = 2100   # L31: 0 -> 2100
$ws.Cells.Item(31, 13).Value = -1212   # M31: None -> -1212
$ws.Cells.Item(31, 14).Value = -2676   # N31: None -> -2676
$ws.Cells.Item(49, 8).Value = 2850   # H49: 3000 -> 2850
$ws.Cells.Item(49, 10).Value = 2850   # J49: 3000 -> 2850
$ws.Cells.Item(49, 12).Value = 8550   # L49: 9000 -> 8550
$ws.Cells.Item(49, 14).Value = -8862   # N49: -9312 -> -8862
$ws.Cells.Item(63, 8).Value = 9304.706   # H63: 4859.25 -> 9304.706
$ws.Cells.Item(63, 9).Value = 4363.3335   # I63: 3976.6667 -> 4363.3335
$ws.Cells.Item(63, 10).Value = 12000   # J63: 7507 -> 12000
$ws.Cells.Item(63, 11).Value = 13090.0005   # K63: 11930.0001 -> 13090.0005
$ws.Cells.Item(63, 12).Value = 36000   # L63: 22521 -> 36000
$ws.Cells.Item(63, 13).Value = -12341.0005   # M63: -11181.0001 -> -12341.0005
$ws.Cells.Item(63, 14).Value = -37498   # N63: -24019 -> -37498
$ws.Cells.Item(66, 8).Value = 9304.706   # H66: 4859.25 -> 9304.706
$ws.Cells.Item(66, 9).Value = 4363.3335   # I66: 3976.6667 -> 4363.3335
$ws.Cells.Item(66, 10).Value = 12000   # J66: 7507 -> 12000
$ws.Cells.Item(66, 11).Value = 39270.0015   # K66: 35790.0003 -> 39270.0015
$ws.Cells.Item(66, 12).Value = 108000   # L66: 67563 -> 108000
$ws.Cells.Item(66, 13).Value = -35526.0015   # M66: -32046.0003 -> -35526.0015
$ws.Cells.Item(66, 14).Value = -115488   # N66: -75051 -> -115488
$ws.Cells.Item(74, 8).Value = 11325   # H74: 9400 -> 11325
$ws.Cells.Item(74, 10).Value = 11325   # J74: 9400 -> 11325
$ws.Cells.Item(74, 12).Value = 33975   # L74: 28200 -> 33975
$ws.Cells.Item(74, 14).Value = -36097   # N74: -30322 -> -36097
$ws.Cells.Item(77, 8).Value = 11325   # H77: 9400 -> 11325
$ws.Cells.Item(77, 10).Value = 11325   # J77: 9400 -> 11325
$ws.Cells.Item(77, 12).Value = 101925   # L77: 84600 -> 101925
$ws.Cells.Item(77, 14).Value = -112533   # N77: -95208 -> -112533
$ws.Cells.Item(93, 8).Value = 3000   # H93: 3095.238 -> 3000
$ws.Cells.Item(93, 9).Value = 1000   # I93: 0 -> 1000
$ws.Cells.Item(93, 11).Value = 3000   # K93: 0 -> 3000
$ws.Cells.Item(93, 13).Value = -1128   # M93: None -> -1128
$ws.Cells.Item(94, 8).Value = 10608.1   # H94: 11342.333 -> 10608.1
$ws.Cells.Item(94, 10).Value = 11453.444   # J94: 12385.125 -> 11453.444
$ws.Cells.Item(94, 12).Value = 34360.33199999999   # L94: 37155.375 -> 34360.33199999999
$ws.Cells.Item(94, 14).Value = -35712.33199999999   # N94: -38507.375 -> -35712.33199999999
$ws.Cells.Item(100, 8).Value = 11411.2   # H100: 11464.728 -> 11411.2
$ws.Cells.Item(100, 9).Value = 1000   # I100: 0 -> 1000
$ws.Cells.Item(100, 10).Value = 12568   # J100: 11464.728 -> 12568
$ws.Cells.Item(100, 11).Value = 3000   # K100: 0 -> 3000
$ws.Cells.Item(100, 12).Value = 37704   # L100: 34394.18399999999 -> 37704
$ws.Cells.Item(100, 13).Value = -2189   # M100: None -> -2189
$ws.Cells.Item(100, 14).Value = -39326   # N100: -36016.18399999999 -> -39326
$ws.Cells.Item(102, 8).Value = 2304.8333   # H102: 3000 -> 2304.8333
$ws.Cells.Item(102, 10).Value = 2304.8333   # J102: 3000 -> 2304.8333
$ws.Cells.Item(102, 12).Value = 6914.499899999999   # L102: 9000 -> 6914.499899999999
$ws.Cells.Item(102, 14).Value = -11782.4999   # N102: -13868 -> -11782.4999
$ws.Cells.Item(103, 8).Value = 5130.2144   # H103: 5766.1816 -> 5130.2144
$ws.Cells.Item(103, 9).Value = 373.75   # I103: 300 -> 373.75
$ws.Cells.Item(103, 10).Value = 7032.8   # J103: 6312.8 -> 7032.8
$ws.Cells.Item(103, 11).Value = 1121.25   # K103: 900 -> 1121.25
$ws.Cells.Item(103, 12).Value = 21098.4   # L103: 18938.4 -> 21098.4
$ws.Cells.Item(103, 13).Value = -242.25   # M103: -21 -> -242.25
$ws.Cells.Item(103, 14).Value = -22856.4   # N103: -20696.4 -> -22856.4
$ws.Cells.Item(114, 8).Value = 3004.4546   # H114: 5050 -> 3004.4546
$ws.Cells.Item(114, 9).Value = 341.5   # I114: 3400 -> 341.5
$ws.Cells.Item(114, 10).Value = 6200   # J114: 10000 -> 6200
$ws.Cells.Item(114, 11).Value = 1024.5   # K114: 10200 -> 1024.5
$ws.Cells.Item(114, 12).Value = 18600   # L114: 30000 -> 18600
$ws.Cells.Item(114, 13).Value = 2229.5   # M114: -6946 -> 2229.5
$ws.Cells.Item(114, 14).Value = -25108   # N114: -36508 -> -25108
$ws.Cells.Item(130, 8).Value = 1754.5333   # H130: 1268.1666 -> 1754.5333
$ws.Cells.Item(130, 10).Value = 3740   # J130: 3800 -> 3740
$ws.Cells.Item(130, 12).Value = 11220   # L130: 11400 -> 11220
$ws.Cells.Item(130, 14).Value = -21260   # N130: -21440 -> -21260
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 18218146   # H70: 18218154 -> 18218146
$ws.Cells.Item(70, 9).Value = 20403820   # I70: 21253824 -> 20403820
$ws.Cells.Item(70, 10).Value = 4200   # J70: 4125 -> 4200
$ws.Cells.Item(70, 11).Value = 20403820   # K70: 21253824 -> 20403820
$ws.Cells.Item(70, 12).Value = 4200   # L70: 4125 -> 4200
$ws.Cells.Item(70, 13).Value = -20403550   # M70: -21253554 -> -20403550
$ws.Cells.Item(70, 14).Value = -4740   # N70: -4665 -> -4740
$ws.Cells.Item(73, 8).Value = 18218146   # H73: 18218154 -> 18218146
$ws.Cells.Item(73, 9).Value = 20403820   # I73: 21253824 -> 20403820
$ws.Cells.Item(73, 10).Value = 4200   # J73: 4125 -> 4200
$ws.Cells.Item(73, 11).Value = 20403820   # K73: 21253824 -> 20403820
$ws.Cells.Item(73, 12).Value = 4200   # L73: 4125 -> 4200
$ws.Cells.Item(73, 13).Value = -20402884   # M73: -21252888 -> -20402884
$ws.Cells.Item(73, 14).Value = -6072   # N73: -5997 -> -6072
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(63, 8).Value = 42923.08   # H63: 53333.332 -> 42923.08
$ws.Cells.Item(63, 10).Value = 42923.08   # J63: 53333.332 -> 42923.08
$ws.Cells.Item(63, 12).Value = 42923.08   # L63: 53333.332 -> 42923.08
$ws.Cells.Item(63, 14).Value = -44421.08   # N63: -54831.332 -> -44421.08
$ws.Cells.Item(64, 8).Value = 0   # H64: 7575 -> 0
$ws.Cells.Item(64, 10).Value = 0   # J64: 7575 -> 0
$ws.Cells.Item(64, 12).Value = 0   # L64: 7575 -> 0
$ws.Cells.Item(64, 14).Value = $null   # N64: clear (was -8025)
$ws.Cells.Item(66, 8).Value = 42923.08   # H66: 53333.332 -> 42923.08
$ws.Cells.Item(66, 10).Value = 42923.08   # J66: 53333.332 -> 42923.08
$ws.Cells.Item(66, 12).Value = 128769.24   # L66: 159999.996 -> 128769.24
$ws.Cells.Item(66, 14).Value = -136257.24   # N66: -167487.996 -> -136257.24
$ws.Cells.Item(67, 8).Value = 0   # H67: 7575 -> 0
$ws.Cells.Item(67, 10).Value = 0   # J67: 7575 -> 0
$ws.Cells.Item(67, 12).Value = 0   # L67: 7575 -> 0
$ws.Cells.Item(67, 14).Value = $null   # N67: clear (was -9135)
$ws.Cells.Item(132, 8).Value = 1913.5758   # H132: 2006.4375 -> 1913.5758
$ws.Cells.Item(132, 9).Value = 2053.889   # I132: 2586.7693 -> 2053.889
$ws.Cells.Item(132, 10).Value = 1745.2   # J132: 1609.3684 -> 1745.2
$ws.Cells.Item(132, 11).Value = 6161.667   # K132: 7760.3079 -> 6161.667
$ws.Cells.Item(132, 12).Value = 5235.6   # L132: 4828.1052 -> 5235.6
$ws.Cells.Item(132, 13).Value = -3631.667   # M132: -5230.3079 -> -3631.667
$ws.Cells.Item(132, 14).Value = -10295.6   # N132: -9888.1052 -> -10295.6
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(25, 8).Value = 15000   # H25: 0 -> 15000
$ws.Cells.Item(25, 10).Value = 15000   # J25: 0 -> 15000
$ws.Cells.Item(25, 12).Value = 15000   # L25: 0 -> 15000
$ws.Cells.Item(25, 14).Value = -15586   # N25: None -> -15586
$ws.Cells.Item(132, 8).Value = 1899.8379   # H132: 1794.5581 -> 1899.8379
$ws.Cells.Item(132, 9).Value = 1337.4482   # I132: 1343.6451 -> 1337.4482
$ws.Cells.Item(132, 10).Value = 3938.5   # J132: 2959.4167 -> 3938.5
$ws.Cells.Item(132, 11).Value = 4012.3446   # K132: 4030.9353 -> 4012.3446
$ws.Cells.Item(132, 12).Value = 11815.5   # L132: 8878.250100000001 -> 11815.5
$ws.Cells.Item(132, 13).Value = -1482.3446   # M132: -1500.9353 -> -1482.3446
$ws.Cells.Item(132, 14).Value = -16875.5   # N132: -13938.2501 -> -16875.5
$ws.Cells.Item(136, 8).Value = 4425.4116   # H136: 4248.028 -> 4425.4116
$ws.Cells.Item(136, 9).Value = 4425.4116   # I136: 4826.9033 -> 4425.4116
$ws.Cells.Item(136, 10).Value = 0   # J136: 659 -> 0
$ws.Cells.Item(136, 11).Value = 13276.2348   # K136: 14480.7099 -> 13276.2348
$ws.Cells.Item(136, 12).Value = 0   # L136: 1977 -> 0
$ws.Cells.Item(136, 13).Value = -10726.2348   # M136: -11930.7099 -> -10726.2348
$ws.Cells.Item(136, 14).Value = $null   # N136: clear (was -7077)
